# Weekly update: insert two new pairs of "Brócoli" price rows (Primera/Segunda)
# into the existing daily-price log, shifting the later rows down.
#
# The new data lands at two insertion points:
#   - before the (old) row 650  -> final rows 650-651
#   - before the (old) row 712  -> final rows 712-713 (after the first insert,
#     that position is still row 712 in the shifted sheet, because the first
#     insert only affected rows at/after 650)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- first insertion point: two new rows at 650-651 ---
$ws.Rows("650:651").Insert()

$ws.Range("A650").Value2 = 8
$ws.Range("B650").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C650").Value2 = "Coquimbo"
$ws.Range("D650").Value2 = 44748
$ws.Range("E650").Value2 = 4
$ws.Range("F650").Value2 = 100112023
$ws.Range("G650").Value2 = "Brócoli"
$ws.Range("H650").Value2 = "Sin especificar"
$ws.Range("I650").Value2 = "Primera"
$ws.Range("J650").Value2 = 2600
$ws.Range("K650").Value2 = 750
$ws.Range("L650").Value2 = 800
$ws.Range("M650").Value2 = 775
$ws.Range("N650").Value2 = "`$/unidad"
$ws.Range("O650").Value2 = "Provincia del Elquí"
$ws.Range("P650").Value2 = 775
$ws.Range("Q650").Value2 = 1
$ws.Range("R650").Value2 = "Hortaliza"

$ws.Range("A651").Value2 = 8
$ws.Range("B651").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C651").Value2 = "Coquimbo"
$ws.Range("D651").Value2 = 44748
$ws.Range("E651").Value2 = 4
$ws.Range("F651").Value2 = 100112023
$ws.Range("G651").Value2 = "Brócoli"
$ws.Range("H651").Value2 = "Sin especificar"
$ws.Range("I651").Value2 = "Segunda"
$ws.Range("J651").Value2 = 1560
$ws.Range("K651").Value2 = 650
$ws.Range("L651").Value2 = 700
$ws.Range("M651").Value2 = 675
$ws.Range("N651").Value2 = "`$/unidad"
$ws.Range("O651").Value2 = "Provincia del Elquí"
$ws.Range("P651").Value2 = 675
$ws.Range("Q651").Value2 = 1
$ws.Range("R651").Value2 = "Hortaliza"

# --- second insertion point: two new rows at (what is now) 712-713 ---
$ws.Rows("712:713").Insert()

$ws.Range("A712").Value2 = 8
$ws.Range("B712").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C712").Value2 = "Coquimbo"
$ws.Range("D712").Value2 = 44747
$ws.Range("E712").Value2 = 4
$ws.Range("F712").Value2 = 100112023
$ws.Range("G712").Value2 = "Brócoli"
$ws.Range("H712").Value2 = "Sin especificar"
$ws.Range("I712").Value2 = "Primera"
$ws.Range("J712").Value2 = 2460
$ws.Range("K712").Value2 = 750
$ws.Range("L712").Value2 = 800
$ws.Range("M712").Value2 = 775
$ws.Range("N712").Value2 = "`$/unidad"
$ws.Range("O712").Value2 = "Provincia del Elquí"
$ws.Range("P712").Value2 = 775
$ws.Range("Q712").Value2 = 1
$ws.Range("R712").Value2 = "Hortaliza"

$ws.Range("A713").Value2 = 8
$ws.Range("B713").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C713").Value2 = "Coquimbo"
$ws.Range("D713").Value2 = 44747
$ws.Range("E713").Value2 = 4
$ws.Range("F713").Value2 = 100112023
$ws.Range("G713").Value2 = "Brócoli"
$ws.Range("H713").Value2 = "Sin especificar"
$ws.Range("I713").Value2 = "Segunda"
$ws.Range("J713").Value2 = 1360
$ws.Range("K713").Value2 = 650
$ws.Range("L713").Value2 = 700
$ws.Range("M713").Value2 = 675
$ws.Range("N713").Value2 = "`$/unidad"
$ws.Range("O713").Value2 = "Provincia del Elquí"
$ws.Range("P713").Value2 = 675
$ws.Range("Q713").Value2 = 1
$ws.Range("R713").Value2 = "Hortaliza"
